$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Last status check on" timestamp in F1
$ws.Range("F1").Value = "Last status check on: 26.01.2022 09:15"

# Row 6 (Shell Olomoucká): refreshed price check
# New current price (B6) / previous price moved to Old Cena (C6)
$ws.Range("B6").Value = 36.9
$ws.Range("C6").Value = 36.7

# Delta is now recorded as a formatted text string, not a number
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "+0.2"
$ws.Range("D6").Style = "Normal"

# Old Datum is now a plain text timestamp instead of a date serial
$ws.Range("E6").Value = "2022-01-26 09:15:12"
$ws.Range("E6").Style = "Normal"
